$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.980.49'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.560.99'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.36'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.782.57'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.561.28'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.11'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.983.23'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.29'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0702'
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.36'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.22'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  -1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.56'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("E27").Value = '  +1.30%  '
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.12'
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.422.82'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +2.85%  '
$ws.Range("E36").Value = '  +8.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.33'
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E39").Value = '  +2.25%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("E43").Value = '  +2.82%  '
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.696.26'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.30'
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0956'
$ws.Range("E51").Value = '  -0.13%  '
